$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update time values for rows 10 and 11 (Item 9 and Item 10)
$ws.Range("C10").Value = 120
$ws.Range("C11").Value = 180

# Recalculate so the SUM formula in C15 reflects the new totals
$excel.Calculate()

# Move/restore the active cell selection to B9
$ws.Range("B9").Select()
